$wb = $excel.ActiveWorkbook
$hungary = $wb.Worksheets.Item("Hungary")

# --- Create "Norway" sheet as a copy of "Hungary", placed right after it ---
$hungary.Copy($null, $hungary)
$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"

# Update the ticket code (B4) before the market name (B2) so the new shared
# strings are appended in the same order as the source data.
$norway.Range("B4").Value = "NGC-2931/T3083/T3069"
$norway.Range("B2").Value = "Norway Market"

# Reset the explicit row heights that were inherited from the Hungary sheet
# back to the sheet's default height.
$norway.Rows.Item(3).AutoFit()
$norway.Rows.Item(4).AutoFit()
$norway.Rows.Item(5).AutoFit()

# Column D on the new sheet is a plain custom width (not an auto best-fit one).
$norway.Columns.Item(4).ColumnWidth = 20.14

# Select the whole sheet (matches the template used for freshly added sheets).
$norway.Cells.Select() | Out-Null

# --- Create "Poland" sheet as a copy of "Hungary", placed right after "Norway" ---
$hungary.Copy($null, $norway)
$poland = $wb.Worksheets.Item($wb.Worksheets.Count)
$poland.Name = "Poland"

$poland.Range("B4").Value = "NGC-2920/T3035/T3118"
$poland.Range("B2").Value = "Poland Market"

$poland.Rows.Item(3).AutoFit()
$poland.Rows.Item(4).AutoFit()
$poland.Rows.Item(5).AutoFit()

$poland.Columns.Item(4).ColumnWidth = 20.14

$poland.Cells.Select() | Out-Null

# "Norway" ends up as the active/selected tab; Select() (rather than
# Activate()) also clears the multi-tab-selection highlight left behind on
# "Hungary" by the two Copy() calls above.
$norway.Select() | Out-Null
